$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.764.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.776.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.13%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "333.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("E7").Value = "  +2.54%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.577"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.12%  "
$ws.Range("E11").Value = "  +5.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.216.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.38%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.779.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.14%  "
$ws.Range("E17").Value = "  +4.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.761.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.37%  "
$ws.Range("E19").Value = "  +10.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "280.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("E30").Value = "  +3.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("E34").Value = "  +1.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.52%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.29"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0359"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "127.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.43%  "
$ws.Range("E44").Value = "  +19.62%  "
$ws.Range("E45").Value = "  +2.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.089.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.59%  "
$ws.Range("E48").Value = "  +4.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "60.80"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.83"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.72%  "
